$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UCROffenseCodeType")
$ws.Activate()

# Insert a new row 58 (shifting 906-999999 block down by one) and give it
# the same "short row" formatting (15pt, custom height) as its neighbours.
$ws.Rows.Item(58).Insert()
$ws.Rows.Item(58).RowHeight = 15

# Re-add the "90I" / "Runaway" Group B Offense code that was dropped earlier.
$ws.Range("A58").Value2 = 909
$ws.Range("B58").Value2 = "90I"
$ws.Range("C58").Value2 = "Runaway"
$ws.Range("D58").Value2 = "90I"
$ws.Range("E58").Value2 = "Runaway"
$ws.Range("F58").Value2 = "Group B"
$ws.Range("G58").Value2 = "Other"
$ws.Range("H58").Value2 = "Group B Offenses (Society)"
$ws.Range("I58").Value2 = "Society"

# Leave the workbook with this sheet active/scrolled/selected, matching the
# state it was saved in.
$ws.Range("D58:E58").Select()
